$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns AD/AE/AF ("Wins"/"Losses"/"Ties"), matching the
# bold/centered/bordered style already used by the rest of row 1 (copy the
# style from AC1, then overwrite the value/text).
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "Wins"

$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "Losses"

$ws.Range("AC1").Copy($ws.Range("AF1"))
$ws.Range("AF1").Value = "Ties"

# Data rows 2-41: team record columns (Wins=74, Losses=88, Ties=0) for
# every player row.
for ($r = 2; $r -le 41; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 88
    $ws.Cells.Item($r, 32).Value = 0
}
